# Applies the MOSFERATU JLCPCB CPL update:
#  - Refresh Mid X / Mid Y / Rotation values for every component row (A2:E63)
#  - Re-apply the Designator sort on Table1 (values already sorted; keeps sortState metadata)
#  - Restore the last-used selection / scroll position on Sheet1

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{ Row = 2; B = -4.973399; C = -31.073399; E = 180 },
    @{ Row = 3; B = -24.956999; C = -18.267599; E = 180 },
    @{ Row = 4; B = -8.874999; C = 8.2074; E = 0 },
    @{ Row = 5; B = 8.791; C = -16.000599; E = 270 },
    @{ Row = 6; B = 25.691; C = -3.343599; E = 90 },
    @{ Row = 7; B = 12.541; C = -15.550599; E = 90 },
    @{ Row = 8; B = -9.424999; C = -15.292599; E = 180 },
    @{ Row = 9; B = 20.041; C = -15.550599; E = 90 },
    @{ Row = 10; B = -25.798999; C = 2.5984; E = 270 },
    @{ Row = 11; B = -25.798999; C = -3.201599; E = 270 },
    @{ Row = 12; B = 12.541; C = -21.350599; E = 90 },
    @{ Row = 13; B = 6.641; C = -2.543599; E = 0 },
    @{ Row = 14; B = -0.474999; C = -8.567599; E = 270 },
    @{ Row = 15; B = -18.706999; C = -29.391599; E = 0 },
    @{ Row = 16; B = -13.874999; C = 3.6574; E = 90 },
    @{ Row = 17; B = 16.291; C = -15.550599; E = 90 },
    @{ Row = 18; B = 7.541; C = -20.550599; E = 180 },
    @{ Row = 19; B = 1.775; C = 14.0234; E = 0 },
    @{ Row = 20; B = -3.624999; C = -15.292599; E = 180 },
    @{ Row = 21; B = 17.3766; C = 4.8234; E = 180 },
    @{ Row = 22; B = 23.1766; C = 4.8234; E = 180 },
    @{ Row = 23; B = 5.041; C = -10.742599; E = 270 },
    @{ Row = 24; B = -18.432999; C = -7.367599; E = 270 },
    @{ Row = 25; B = 0.325; C = 1.6904; E = 180 },
    @{ Row = 26; B = 21.941; C = -3.343599; E = 90 },
    @{ Row = 27; B = 1.075; C = -24.258599; E = 0 },
    @{ Row = 28; B = 10.058953; C = 21.571; E = 90 },
    @{ Row = 29; B = 19.791; C = -23.400599; E = 0 },
    @{ Row = 30; B = 8.891; C = 4.0564; E = 180 },
    @{ Row = 31; B = 14.691; C = -5.393599; E = 180 },
    @{ Row = 32; B = -14.798999; C = -22.733599; E = 0 },
    @{ Row = 33; B = 0.0056; C = -31.827399; E = 90 },
    @{ Row = 34; B = 0.279; C = 5.4194; E = 0 },
    @{ Row = 35; B = -8.828999; C = 11.9364; E = 0 },
    @{ Row = 36; B = -24.460999; C = -21.996599; E = 0 },
    @{ Row = 37; B = -18.378999; C = 4.4784; E = 180 },
    @{ Row = 38; B = 23.77; C = -16.046599; E = 90 },
    @{ Row = 39; B = 13.012; C = 11.534; E = 90 },
    @{ Row = 40; B = 23.77; C = -10.788599; E = 270 },
    @{ Row = 41; B = -8.828999; C = 15.6444; E = 180 },
    @{ Row = 42; B = -14.728999; C = -15.271599; E = 0 },
    @{ Row = 43; B = -25.777999; C = -8.505599; E = 90 },
    @{ Row = 44; B = -7.107999; C = 20.746; E = 0 },
    @{ Row = 45; B = -0.495999; C = -2.813599; E = 270 },
    @{ Row = 46; B = -25.777999; C = -13.763599; E = 270 },
    @{ Row = 47; B = 11.787; C = -26.329599; E = 0 },
    @{ Row = 48; B = 6.379; C = -24.279599; E = 180 },
    @{ Row = 49; B = -24.460999; C = -25.704599; E = 0 },
    @{ Row = 50; B = -11.590999; C = 21.521; E = 270 },
    @{ Row = 51; B = -10.815999; C = 26.004; E = 0 },
    @{ Row = 52; B = -5.220999; C = -22.729599; E = 180 },
    @{ Row = 53; B = 5.062; C = -16.046599; E = 270 },
    @{ Row = 54; B = -5.220999; C = -19.021599; E = 180 },
    @{ Row = 55; B = 1.354; C = -16.046599; E = 270 },
    @{ Row = 56; B = 0.579; C = -20.529599; E = 180 },
    @{ Row = 57; B = -5.220999; C = -26.437599; E = 180 },
    @{ Row = 58; B = -21.294999; C = 0.1614; E = 0 },
    @{ Row = 59; B = 9.304; C = 11.534; E = 270 },
    @{ Row = 60; B = -24.460999; C = -29.412599; E = 0 },
    @{ Row = 61; B = 10.079; C = 16.017; E = 180 },
    @{ Row = 62; B = -7.174999; C = 0.5574; E = 90 },
    @{ Row = 63; B = -7.224999; C = -9.542599; E = 270 }
)

foreach ($item in $data) {
    $ws.Cells.Item($item.Row, 2).Value = $item.B
    $ws.Cells.Item($item.Row, 3).Value = $item.C
    $ws.Cells.Item($item.Row, 5).Value = $item.E
}

# Re-sort the table by Designator (column A) - matches the sortState recorded in table1.xml
try {
    $tbl = $ws.ListObjects.Item(1)
    $tbl.Sort.SortFields.Clear()
    $tbl.Sort.SortFields.Add($ws.Range("A2:A63"))
    $tbl.Sort.Header = 1
    $tbl.Sort.Apply()
} catch {
    # Sort surface may be unavailable; data is already in sorted order so this is cosmetic.
}

# Restore scroll position / selection as last saved
try {
    $excel.ActiveWindow.ScrollRow = 49
    $excel.ActiveWindow.ScrollColumn = 1
} catch {
}
$ws.Range("A60").Select()
